$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 10 new rows above existing row 2, pushing current data rows (2-8) down to (12-18)
$ws.Range("A2:A11").EntireRow.Insert()

# The inserted rows pick up the header row's formatting by default; clear it
# so the new data rows stay unformatted like the original data rows.
$ws.Range("A2:C11").ClearFormats()

# New data to place into rows 2-11
$newData = @(
    @(1.223868489265442, 1.015716195106506, 0.2379320114850998),
    @(-0.2808452844619751, -0.5580254197120667, 0.9086620807647704),
    @(-0.4292855560779571, -1.047786593437195, 0.4741841554641723),
    @(0.2979495227336883, -0.4011857509613037, 0.1365283876657486),
    @(-0.0152716310694813, -0.2217440903186798, -0.104763388633728),
    @(-0.0087048299610614, 0.0255036242306232, -0.0922406539320945),
    @(0.0565050356090068, -0.0242818929255008, 0.0531452745199203),
    @(0.0300851128995418, 0.0068722339347004, -0.0010690141934901),
    @(0.1162171140313148, -0.3998112976551056, 0.1134682223200798),
    @(-0.1954768747091293, -0.8119926452636719, 0.1353066563606262)
)

$row = 2
foreach ($values in $newData) {
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $row++
}
